$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.841.84"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "'2.306.65"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'317.20"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'104.13"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'39.97"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'8.53"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").Value = "'15.42"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'2.657.79"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'2.306.02"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'42.781.87"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'14.44"
$ws.Range("E20").Value = "  +34.87%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'74.11"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "'3.54"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").Value = "'267.40"
$ws.Range("E24").Value = "  -5.29%  "
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'10.99"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "'6.73"
$ws.Range("E29").Value = "  +13.74%  "
$ws.Range("D30").Value = "'22.67"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'37.53"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").Value = "'165.88"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'0.0886"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").Value = "'2.61"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").Value = "'4.58"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "'0.0355"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").Value = "'1.59"
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("D42").Value = "'70.53"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "'96.06"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").Value = "'82.38"
$ws.Range("E47").Value = "  +4.45%  "
$ws.Range("D48").Value = "'115.20"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").Value = "'1.682.89"
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("D50").Value = "'8.85"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").Value = "'5.21"
$ws.Range("E51").Value = "  -2.23%  "
